# The document contains a paragraph whose content is a single Word field:
#   { m:'doc.html'.fromHTMLURI() }
# built from w:fldChar (begin/end) + w:instrText runs. The commit switches
# the parser to a literal-token reader (TokenIteratorFieldRewriterSplit),
# so the field is rewritten as plain text runs (w:t) spelling out the same
# token stream, without surrounding spaces, wrapped in literal "{" / "}"
# runs - while keeping the _GoBack bookmark untouched in the middle.

$d = $word.ActiveDocument

# Locate the paragraph that holds the field (the one with fldChar/instrText).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
    }
}

$apos = "'"

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r><w:t>{</w:t></w:r>' + `
    '<w:r><w:t>m</w:t></w:r>' + `
    '<w:r><w:t>:</w:t></w:r>' + `
    '<w:r><w:t>' + $apos + '</w:t></w:r>' + `
    '<w:r><w:t>doc.html</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>' + $apos + '.fromHTMLURI()</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' + `
    '</w:p>'

# Replace that paragraph's content (InsertXML replaces the contents of the
# exact range it is called on) with the literal-text run sequence.
$target.Range.InsertXML($newXml)
